$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row 13 cell text
$ws.Range("B13").Value2 = "Buchungsdatum"
$ws.Range("D13").Value2 = "Betrag"
$ws.Range("J13").Value2 = "Mögliche Zahlpläne"

# Row 13 height change
$ws.Rows.Item(13).RowHeight = 13.8

# Update selection to A14
$ws.Range("A14").Select() | Out-Null

# Adjust header/footer margins (cosmetic, matches re-save rounding)
$ws.PageSetup.HeaderMargin = 36.8503937007874
$ws.PageSetup.FooterMargin = 36.8503937007874
